# Apply cryptos.xlsx price/volume update ("Updated cryptos list" GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as literal text (e.g. "215.00",
# "8.50", thousand-grouped "90.424.79"). Flip the numeric-looking cells to a Text
# number format before writing so Excel does not coerce them to real numbers and
# strip significant trailing/grouping digits, then restore the default style so no
# stray formatting is left behind.
$numericPriceCells = "D5","D6","D7","D8","D11","D13","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D28","D30","D31","D32","D33","D34","D35","D36","D38","D39","D41","D44","D45","D46","D47","D50","D51"
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '90.424.79'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = '3.186.76'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '215.00'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").Value = '620.12'
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("D7").Value = '0.398'
$ws.Range("E7").Value = '  +3.48%  '
$ws.Range("D8").Value = '0.691'
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '3.183.18'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").Value = '0.578'
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("E12").Value = '  -6.01%  '
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  -3.67%  '
$ws.Range("D14").Value = '90.123.40'
$ws.Range("E14").Value = '  +2.58%  '
$ws.Range("D15").Value = '3.778.27'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").Value = '33.07'
$ws.Range("E16").Value = '  -3.13%  '
$ws.Range("D17").Value = '5.26'
$ws.Range("E17").Value = '  -4.04%  '
$ws.Range("D18").Value = '3.181.36'
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("D19").Value = '3.29'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").Value = '0.0000207'
$ws.Range("E20").Value = '  +44.74%  '
$ws.Range("D21").Value = '13.45'
$ws.Range("E21").Value = '  -4.01%  '
$ws.Range("D22").Value = '438.52'
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").Value = '8.62'
$ws.Range("E23").Value = '  -4.13%  '
$ws.Range("D24").Value = '5.07'
$ws.Range("E24").Value = '  -5.09%  '
$ws.Range("D25").Value = '5.16'
$ws.Range("E25").Value = '  -4.13%  '
$ws.Range("D26").Value = '11.68'
$ws.Range("E26").Value = '  -6.09%  '
$ws.Range("D27").Value = '3.356.06'
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("D28").Value = '75.42'
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '0.172'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '4.17'
$ws.Range("E32").Value = '  +26.43%  '
$ws.Range("D33").Value = '8.50'
$ws.Range("E33").Value = '  -4.07%  '
$ws.Range("D34").Value = '537.02'
$ws.Range("E34").Value = '  -5.96%  '
$ws.Range("D35").Value = '7.01'
$ws.Range("E35").Value = '  -3.38%  '
$ws.Range("D36").Value = '1.87'
$ws.Range("E36").Value = '  -4.77%  '
$ws.Range("E37").Value = '  -8.69%  '
$ws.Range("D38").Value = '22.10'
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("D39").Value = '22.35'
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("E40").Value = '  -8.12%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -4.08%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '0.375'
$ws.Range("E44").Value = '  -6.64%  '
$ws.Range("D45").Value = '150.44'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '43.63'
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '172.74'
$ws.Range("E47").Value = '  -4.10%  '
$ws.Range("E48").Value = '  -8.28%  '
$ws.Range("E49").Value = '  -8.46%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = '0.611'
$ws.Range("E50").Value = '  -3.57%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '4.06'
$ws.Range("E51").Value = '  -4.33%  '

foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}

